$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update roll_no and student_name (shared strings) for row 2
$ws.Range("B2").Value = "20NU1A0503"
$ws.Range("C2").Value = "ATTA HARIKA"

# Update numeric data for row 2
$ws.Range("A2").Value = 3
$ws.Range("D2").Value = 8.9499999999999993
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 82.01

# Move active selection to E2
$ws.Range("E2").Select()
